$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 202 (pushes existing rows 202-238 down to 203-239)
$ws.Rows.Item(202).Insert()

# Populate the new row 202 with the new weekly price record
$ws.Cells.Item(202, 1).Value = 10
$ws.Cells.Item(202, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(202, 3).Value = "La Araucanía"
$ws.Cells.Item(202, 4).Value = (Get-Date -Year 2023 -Month 3 -Day 10 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(202, 5).Value = 9
$ws.Cells.Item(202, 6).Value = "Fruta"
$ws.Cells.Item(202, 7).Value = 100104
$ws.Cells.Item(202, 8).Value = "Frutos de pepita"
$ws.Cells.Item(202, 9).Value = 100104003
$ws.Cells.Item(202, 10).Value = "Membrillo"
$ws.Cells.Item(202, 11).Value = "Champion"
$ws.Cells.Item(202, 12).Value = "Primera"
$ws.Cells.Item(202, 13).Value = 50
$ws.Cells.Item(202, 14).Value = 14000
$ws.Cells.Item(202, 15).Value = 14000
$ws.Cells.Item(202, 16).Value = 14000
$ws.Cells.Item(202, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(202, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(202, 19).Value = 778
$ws.Cells.Item(202, 20).Value = 18
